# ContactPage.xlsx — add a "Locator Type" column (D) classifying each
# locator as Xpath (row 2) or CSS (rows 3-28), matching the mobile
# Android locator work described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: header + values ------------------------------------
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D2").Value = "Xpath"
for ($r = 3; $r -le 28; $r++) {
    $ws.Range("D$r").Value = "CSS"
}

# --- Header cell styling (bold, 13pt, Helvetica Neue) ------------------
# Order matters for minimal style-table churn: Bold first (reuses the
# existing bold font), then Size, then Name.
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 13
$ws.Range("D1").Font.Name = "Helvetica Neue"

# --- Row 1 height & column D width --------------------------------------
$ws.Rows.Item(1).RowHeight = 17
$ws.Columns.Item(4).ColumnWidth = 26.666666666666668

# --- View state: scroll so row 16 / col B is top-left, select C28 ------
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 2
$ws.Range("C28").Select() | Out-Null
